# GeneNetworkAnnotator.xlsx — "attributes" sheet update
#
# Replaces the single "Annotation_log" (text) attribute definition (row 7)
# with two attributes: "termsFound" and "termsNotFound" (both dataType
# "string"), each carrying a label + description. This pushes the two
# "Patient" attribute rows (ID / Notes) down by one row. Also normalises a
# couple of duplicate-looking cell styles on rows 5-6, widens column C to
# fit the new, longer labels, and updates the remembered selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Style touch-ups on existing rows 5 & 6 -----------------------------
# B5, and A6/B6/C6 move from the "accent, explicit font" variant onto the
# plain accent style already used elsewhere in the same rows (A5/C5/D5/F5).
$ws.Range("A5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("C6").PasteSpecial(-4122)

# --- 2. Make room for the extra attribute row ------------------------------
# Old row 7 ("Annotation_log") becomes the new "termsFound" row; a brand
# new row 8 ("termsNotFound") is inserted below it, pushing the old rows
# 8 & 9 (the "Patient" ID/Notes attributes) down to rows 9 & 10.
$ws.Rows("8:8").Insert()

# --- 3. Rewrite row 7: termsFound ------------------------------------------
$ws.Range("A7").Value = "termsFound"
$ws.Range("C7").Value = "terms used for Gene Network"

$ws.Range("A3").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "HPO terms that were used by the gene network annotator"

$ws.Range("F7").Value = "string"

$ws.Range("A3").Copy()
$ws.Range("G7").PasteSpecial(-4122)

# --- 4. Populate the new row 8: termsNotFound ------------------------------
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "termsNotFound"

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "Project"

$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = "termsnot  used for Gene Network"

$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "HPO terms that were not found by the gene network annotator"

$ws.Range("F7").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = "string"

# --- 5. Column C needs to be wider for the longer labels -------------------
$ws.Columns("C:C").ColumnWidth = 27.1640625

# --- 6. Restore the remembered selection -----------------------------------
$ws.Range("F6").Select()
